# The workbook tracks weekly "Poroto verde" (green bean) price quotes from
# "Comercializadora del Agro de Limarí". A new week's worth of data (two
# rows: one for the "Magnum" variety and one for "Sin especificar") is
# inserted right before the existing row 99, pushing all subsequent rows
# down by two. The sheet's used range grows from A1:R180 to A1:R182.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at positions 99-100 (existing rows 99.. shift to 101..).
$ws.Range("A99:A100").EntireRow.Insert()

# --- New row 99: Magnum, fecha 44651 ---
$ws.Cells.Item(99, 1).Value = 2
$ws.Cells.Item(99, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(99, 3).Value = "Coquimbo"
$ws.Cells.Item(99, 4).Value = 44651
$ws.Cells.Item(99, 5).Value = 4
$ws.Cells.Item(99, 6).Value = 100112031
$ws.Cells.Item(99, 7).Value = "Poroto verde"
$ws.Cells.Item(99, 8).Value = "Magnum"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 400
$ws.Cells.Item(99, 11).Value = 16000
$ws.Cells.Item(99, 12).Value = 18000
$ws.Cells.Item(99, 13).Value = 17000
$ws.Cells.Item(99, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(99, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(99, 16).Value = 680
$ws.Cells.Item(99, 17).Value = 25
$ws.Cells.Item(99, 18).Value = "Hortaliza"

# --- New row 100: Sin especificar, fecha 44651 ---
$ws.Cells.Item(100, 1).Value = 2
$ws.Cells.Item(100, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(100, 3).Value = "Coquimbo"
$ws.Cells.Item(100, 4).Value = 44651
$ws.Cells.Item(100, 5).Value = 4
$ws.Cells.Item(100, 6).Value = 100112031
$ws.Cells.Item(100, 7).Value = "Poroto verde"
$ws.Cells.Item(100, 8).Value = "Sin especificar"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 400
$ws.Cells.Item(100, 11).Value = 23000
$ws.Cells.Item(100, 12).Value = 25000
$ws.Cells.Item(100, 13).Value = 24000
$ws.Cells.Item(100, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(100, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(100, 16).Value = 960
$ws.Cells.Item(100, 17).Value = 25
$ws.Cells.Item(100, 18).Value = "Hortaliza"

# Apply the same date-serial number format used by the rest of column D
# to the two new D cells so they render as dates like their neighbours.
$ws.Range("D99:D100").NumberFormat = "YYYY-MM-DD HH:MM:SS"
